$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: CheckIn ---
$ws.Range("A2").Value = "CheckIn"
$ws.Range("B2").Value = '[0,0,0,0,0,0,0,0,0,0,0,0,0,0,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,0,0,0,0,0,0,0,0,0,0,0,0,0,0]'
$ws.Range("C2").Value = 6
$ws.Range("D2").Value = '["Business", "Speakers", "Logistics"]'

# --- Row 3: Auditorio ---
$ws.Range("A3").Value = "Auditorio"
$ws.Range("B3").Value = '[0,0,0,0,0,0,0,0,0,0,0,0,0,0,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,0,0,0,0,0,0,0,0,0,0,0,0,0,0]'
$ws.Range("C3").Value = 2
$ws.Range("D3").Value = '["Logistics"]'

# --- Row 4: Refeicoes ---
$ws.Range("A4").Value = "Refeicoes"
$ws.Range("B4").Value = '[0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,1,1,0,0,0,0,0,1,1,0,0,0,0,0,1,1,0,0,0,0,0,1,1,0,0,0,0,0,1,1,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0]'
$ws.Range("C4").Value = 6
$ws.Range("D4").Value = '["Business", "Logistics"]'

# --- Row 5: Divulgacao ---
$ws.Range("A5").Value = "Divulgacao"
$ws.Range("B5").Value = '[0,0,0,0,0,0,0,0,0,0,0,0,0,0,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,0,0,0,0,0,0,0,0,0,0,0,0,0,0]'
$ws.Range("C5").Value = 5
$ws.Range("D5").Value = '["Marketing", "Volunteer"]'

# --- Row 6: Workshops ---
$ws.Range("A6").Value = "Workshops"
$ws.Range("B6").Value = '[0,0,0,0,0,0,0,0,0,0,0,0,0,0,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,0,0,0,0,0,0,0,0,0,0,0,0,0,0]'
$ws.Range("C6").Value = 4
$ws.Range("D6").Value = '[]'

# --- Row 7: MontagemDesmontagem ---
$ws.Range("A7").Value = "MontagemDesmontagem"
$ws.Range("B7").Value = '[1,1,1,1,1,1,1,1,1,1,1,1,1,1,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,1,1,1,1,1,1,1,1,1,1,1,1,1,1]'
$ws.Range("C7").Value = 6
$ws.Range("D7").Value = '[]'

# --- Row 8: CoffeeBreak ---
$ws.Range("A8").Value = "CoffeeBreak"
$ws.Range("B8").Value = '[0,0,0,0,0,0,0,0,0,0,0,0,0,0,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,0,0,0,0,0,0,0,0,0,0,0,0,0,0]'
$ws.Range("C8").Value = 3
$ws.Range("D8").Value = '[]'

# --- Row 9: MarketingTurno ---
$ws.Range("A9").Value = "MarketingTurno"
$ws.Range("B9").Value = '[0,0,0,0,0,0,0,0,0,0,0,0,0,0,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,0,0,0,0,0,0,0,0,0,0,0,0,0,0]'
$ws.Range("C9").Value = 2
$ws.Range("D9").Value = '["Marketing"]'

# --- Column A width ---
$ws.Columns.Item(1).ColumnWidth = 21.1

# --- Selection ---
$ws.Range("D10").Select()
